# Applies the "Removing Isolation code and updating files" edit to the
# DataInsightsByMAQSoftwareChecklist workbook.
#
# Summary of the change:
#  - The "Isolation" BVT feature (old row 14, S no 9) is removed entirely.
#  - "Reset" (old row 15 / S no 10) becomes S no 9; its Steps/Output text is
#    rewritten so it no longer references the removed isolation feature.
#  - "Undo" (old rows 16-17 / S no 11) becomes S no 10.
#  - "Presentation mode" (old row 18 / S no 12) becomes S no 11.
#  - Column E is widened and the saved selection / row heights change too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVT")

# ---------------------------------------------------------------------
# Row 14: used to be the "Isolation" feature row (S no 9). All of its
# content goes away; only C14/D14/E14 survive as blank, wrap-styled cells.
# ---------------------------------------------------------------------
$ws.Range("A14").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Rows.Item(14).AutoFit()

# ---------------------------------------------------------------------
# Row 15: was "Reset" with S no 10 -> now S no 9. Steps/Output rewritten.
# ---------------------------------------------------------------------
$ws.Range("A15").Value = 9
$ws.Range("D15").Value = "1. Click on any bar/column/brick" + [char]10 + "2. Click on Reset button"
$ws.Range("E15").Value = "The selected bar/brick/column gets higher opacity and others get lower opacity" + [char]10 + "Chart is redrawn and opacity gets reset to normal"
$ws.Rows.Item(15).RowHeight = 45

# ---------------------------------------------------------------------
# Row 16/17: was "Undo" with S no 11 -> now S no 10. Text unchanged.
# ---------------------------------------------------------------------
$ws.Range("A16").Value = 10

# ---------------------------------------------------------------------
# Row 18: was "Presentation mode" with S no 12 -> now S no 11.
# ---------------------------------------------------------------------
$ws.Range("A18").Value = 11

# ---------------------------------------------------------------------
# Column E is widened.
# ---------------------------------------------------------------------
$ws.Range("E1").ColumnWidth = 66.3

# ---------------------------------------------------------------------
# Saved cursor/selection moves to E24.
# ---------------------------------------------------------------------
$ws.Range("E24").Select()
